$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# Row 21: Tom D'Angelo
$ws.Range("A21").Value = "Tom D'Angelo"
foreach ($col in @("C","D","E","I","K","N","O","Q","T","V")) {
    $ws.Range($col + "21").Value = "x"
}
$ws.Range("AK21").Value = 10
$ws.Range("AL21").Value = "DM"
$ws.Range("AM20").Copy()
$ws.Range("AM21").PasteSpecial(-4122)
$ws.Range("AM21").Value = 43441

# Row 22: Bill Plunkett
$ws.Range("A22").Value = "Bill Plunkett"
foreach ($col in @("E","F","I","O","R","V")) {
    $ws.Range($col + "22").Value = "x"
}
$ws.Range("AK22").Value = 6
$ws.Range("AL22").Value = "Email"
$ws.Range("AM20").Copy()
$ws.Range("AM22").PasteSpecial(-4122)
$ws.Range("AM22").Value = 43441

$ws.Range("A22").Select()
